# Adds the "averaged_world_loocv" worksheet (a per-country LOOCV-averaged
# categorical clustering summary) after the existing "raw33d_dissim" sheet,
# fills in its data, and updates the selection/active-sheet state on both
# sheets to match the saved workbook state.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet, name + position it after the existing sheet ---
$ws2 = $wb.Worksheets.Add()
$ws2.Name = "averaged_world_loocv"
$ws2.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

# Re-fetch worksheet handles by name (Move() invalidates earlier references)
$ws1 = $wb.Worksheets.Item("raw33d_dissim")
$ws2 = $wb.Worksheets.Item("averaged_world_loocv")

# --- Populate averaged_world_loocv with the per-country cluster summary ---
$data = @(
  @("cluster_id","USA","UK","South_africa","Australia","Germany","Japan","CHN","HK","France","Chile","Mexico","Spain","Portugal","Brazil","Israel","Russia","Egypt","Qatar","India"),
  @(0,"Private","Hostile","Private","Private","Private","Private","Public","Public","Private","Private","Private","Private","Private","Private","Private","Private","Private","Private","Private"),
  @(1,"Public","Private","Public","Public","Public","Public","Hostile","Private","Public","Public","Public","Public","Public","Public","Public","Public","Public","Public","Public"),
  @(2,"Hostile","Public","Hostile","Hostile","Hostile","Hostile","Private","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile","Hostile")
)

for ($r = 0; $r -lt $data.Count; $r++) {
    for ($c = 0; $c -lt $data[$r].Count; $c++) {
        $ws2.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# --- Restore the original sheet's selection (no longer the active tab) ---
$ws1.Activate()
$ws1.Range("A2:B4").Select()

# --- The new sheet ends up active, with its own lingering selection ---
$ws2.Activate()
$ws2.Range("U10").Select()
